$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the CasesTab query in B2: drop the trailing Cohort column
#    (the stashed-changes commit removes `co.cohort_description` AS `Cohort`)
$newB2 = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Bluetick Hound','Welsh Springer Spaniel','Wheaten Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in ['Bladder, Urethra']`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newB2

# 2. Row heights shrank slightly (re-wrap after the text edit / font-metric change)
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 261
$ws.Rows.Item(4).RowHeight = 261

# 3. View: scrolled one row up, and active/selected cell moved to B2
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select() | Out-Null
